$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of variable names for rows 3-46 (row 2 "aaq_dum" stays put).
# This reflects the variable list being re-sorted/re-specified, shifting each
# subsequent row's label "up" by one, and dropping the former rows 47-62.
$names = @(
    "aca_impa",
    "activ_yn",
    "age",
    "anx_score",
    "assault_sex",
    "belong1",
    "binge_fr",
    "brs_dum",
    "dep_impa",
    "dep_secret",
    "deprawsc",
    "discrim",
    "divers",
    "drug_mar",
    "drugs_yn",
    "dx_adhd",
    "dx_anx",
    "dx_bi",
    "dx_dep",
    "dx_pers",
    "dx_tr",
    "ed_any",
    "fincur",
    "flourish_dum",
    "gad7_impa",
    "gender_noncis",
    "gpa_sr",
    "inf",
    "ins_cover",
    "meds_count",
    "meds_sti",
    "mh_stigma",
    "military",
    "percneed_cur",
    "persist",
    "psyhx",
    "race",
    "religios",
    "satisfied_overall",
    "school2_type",
    "sexual",
    "sib_freq",
    "talk",
    "trauma_year"
)

$startRow = 3
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

# Remove the now-unused trailing rows (former rows 47-62) entirely, along
# with their B/C values, so the sheet's used range shrinks to A1:C46.
$lastRow = 62
$ws.Range("A47:C$lastRow").ClearContents() | Out-Null
